{"js": "// Add \"S3\" to the AWS tools list in the \"Technologies:\" bullet.\n// Before: \"...autoscaling, EC2) \"\n// After:  \"...autoscaling, EC2, S3) \"\nconst body = context.document.body;\n\n// Locate the unique run of text ending the AWS parenthetical list.\nconst results = body.search(\"autoscaling, EC2)\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find target text \"autoscaling, EC2)\" to update.');\n}\n\n// Replace it with the same text plus the new \", S3\" entry, preserving\n// surrounding formatting by replacing in place (same run/paragraph).\nresults.items[0].insertText(\"autoscaling, EC2, S3)\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Add \"S3\" to the AWS tools list in the \"Technologies:\" bullet.\n# Before: \"...autoscaling, EC2) \"\n# After:  \"...autoscaling, EC2, S3) \"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"autoscaling, EC2)\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"autoscaling, EC2, S3)\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceOne = 1 -> replace just the single (unique) match.\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 1) | Out-Null\n"}
